
# Helper: find a paragraph index whose Range.Text contains the given
# (case-sensitive) substring. Using .Contains() because the -like/-match
# operators in this host are not reliably case-sensitive.
function Find-ParaIndex($d, $matchText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.Contains($matchText)) { return $i }
    }
    return -1
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Edit 1: table cell that used to read "Overfitting " becomes an
# empty (bold, right/centre-tab-enabled) paragraph.
# ---------------------------------------------------------------
$idx = Find-ParaIndex $d "Overfitting "
if ($idx -eq -1) { throw "Could not locate the 'Overfitting ' table paragraph" }
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="782BB64C" w14:textId="3D8F1F18" w:rsidR="00F16D1A" w:rsidRPr="007162DC" w:rsidRDefault="00F16D1A" w:rsidP="00C50891"><w:pPr><w:tabs><w:tab w:val="center" w:pos="1888"/></w:tabs><w:jc w:val="both"/><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# ---------------------------------------------------------------
# Edit 2: insert three new paragraphs right after "Learning rate
# gradually decreasing after some Epochs" (and before the two blank
# paragraphs that precede "Further suggestions for improvement:").
# ---------------------------------------------------------------
$idx = Find-ParaIndex $d "Learning rate gradually decreasing after some Epochs"
if ($idx -eq -1) { throw "Could not locate the 'Learning rate...' paragraph" }
$p = $d.Paragraphs.Item($idx + 1)
$r = $p.Range
$r.Collapse(1)
$r.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>But in terms of lighter model Conv3D Model 1 was chosen since It has lighter parameter and accuracy is next better to CNN+LSTM with GRU</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>(Training Accuracy: 9</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>%, Validation Accuracy: 8</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>%)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# ---------------------------------------------------------------
# Edit 3: add a lastRenderedPageBreak marker before "Further
# suggestions for improvement:".
# ---------------------------------------------------------------
$idx = Find-ParaIndex $d "Further suggestions for improvement:"
if ($idx -eq -1) { throw "Could not locate the 'Further suggestions...' paragraph" }
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="28E43CB3" w14:textId="057A73BA" w:rsidR="002023F5" w:rsidRPr="002023F5" w:rsidRDefault="003232AB" w:rsidP="002023F5"><w:pPr><w:pStyle w:val="Heading1"/><w:shd w:val="clear" w:color="auto" w:fill="F5F5F5"/><w:spacing w:before="150" w:beforeAutospacing="0" w:after="225" w:afterAutospacing="0"/><w:rPr><w:color w:val="333333"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00E61F75"><w:rPr><w:color w:val="333333"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t>Further suggestions for improvement:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# ---------------------------------------------------------------
# Edit 4: merge the two runs split by a lastRenderedPageBreak around
# "for sequence information before / finally passing it ..." into a
# single run (and drop the page-break marker).
# ---------------------------------------------------------------
$idx = Find-ParaIndex $d "for sequence information before"
if ($idx -eq -1) { throw "Could not locate the 'for sequence information before' paragraph" }
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="5FB49E70" w14:textId="4296A928" w:rsidR="003232AB" w:rsidRPr="00E61F75" w:rsidRDefault="00116A67" w:rsidP="00E61F75"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/></w:rPr><w:t>Using Transfer L</w:t></w:r><w:r w:rsidRPr="00E61F75"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/></w:rPr><w:t>earning</w:t></w:r><w:r w:rsidRPr="00E61F75"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>:</w:t></w:r><w:r w:rsidR="003232AB" w:rsidRPr="00E61F75"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> Using a pre-trained </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/></w:rPr><w:t>ResN</w:t></w:r><w:r w:rsidR="003232AB" w:rsidRPr="00116A67"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/></w:rPr><w:t>et50/</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/></w:rPr><w:t>ResNet152/Inception V3</w:t></w:r><w:r w:rsidR="003232AB" w:rsidRPr="00E61F75"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> to identify the initial feature vectors and passing them further to a </w:t></w:r><w:r w:rsidR="003232AB" w:rsidRPr="00116A67"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/></w:rPr><w:t>RNN</w:t></w:r><w:r w:rsidR="003232AB" w:rsidRPr="00E61F75"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> for sequence information before finally passing it to a softmax layer for classification of gestures. (This was attempted but </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>other pre-trained models couldn’t be tested</w:t></w:r><w:r w:rsidR="003232AB" w:rsidRPr="00E61F75"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> due to lack of time</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> and disk space in the </w:t></w:r><w:r w:rsidRPr="00116A67"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">nimblebox.ai </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>platform</w:t></w:r><w:r w:rsidR="003232AB" w:rsidRPr="00E61F75"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>.)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

Write-Output "All edits applied"
